# Update symbol list (cryptos.xlsx) - Fri Dec 16 19:49:00 UTC 2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (to preserve exact
# formatting such as trailing/leading zeros). Force the number format to
# "Text" on every D cell we touch before writing so Excel doesn't silently
# convert the string into a floating point number.
$priceCells = @(
    "D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14",
    "D15","D16","D18","D19","D20","D21","D24","D26","D27","D40","D41",
    "D42","D43","D44","D46","D48","D49","D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price ("D") column updates ---
$ws.Range("D2").Value  = "243.24"
$ws.Range("D3").Value  = "23.34"
$ws.Range("D4").Value  = "5.646"
$ws.Range("D5").Value  = "0.05805"
$ws.Range("D6").Value  = "3.409"
$ws.Range("D7").Value  = "6.468"
$ws.Range("D8").Value  = "1.317"
$ws.Range("D9").Value  = "0.7976"
$ws.Range("D10").Value = "0.1459"
$ws.Range("D11").Value = "0.07626"
$ws.Range("D12").Value = "0.03217"
$ws.Range("D13").Value = "0.02959"
$ws.Range("D14").Value = "0.09243"
$ws.Range("D15").Value = "0.001679"
$ws.Range("D16").Value = "3.323"
$ws.Range("D18").Value = "0.0005991"
$ws.Range("D19").Value = "0.006225"
$ws.Range("D20").Value = "0.005462"
$ws.Range("D21").Value = "0.001068"
$ws.Range("D24").Value = "2.194"
$ws.Range("D26").Value = "0.1240"
$ws.Range("D27").Value = "0.0009998"
$ws.Range("D40").Value = "0.04277"
$ws.Range("D41").Value = "0.007138"
$ws.Range("D44").Value = "0.009533"
$ws.Range("D46").Value = "0.00005433"
$ws.Range("D48").Value = "0.7853"
$ws.Range("D49").Value = "0.1046"
$ws.Range("D50").Value = "0.00002100"

# --- Row 18 ("One" / ONE) Volume column label gained a "Worstin24h" suffix ---
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Row 45 (ACDXExchange) lost its "Worstin24h" suffix ---
$ws.Range("E45").Value = "44ACDXExchangeACXT"

# --- Rows 42 and 43 swapped: BKEXToken <-> CEJI (with refreshed prices) ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.003600"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1054"
$ws.Range("E43").Value = "42BKEXTokenBKK"
